$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,15

$arr[0,0] = 0
$arr[0,1] = 0.4667726049514727
$arr[0,2] = 0.1115831045428308
$arr[0,3] = 0
$arr[0,4] = 0.1118171829289238
$arr[0,5] = 0.4443680307746263
$arr[0,6] = 0.00246826798140105
$arr[0,7] = 0
$arr[0,8] = 0.8886284222073222
$arr[0,9] = 0
$arr[0,10] = 0.2659201194302483
$arr[0,11] = 0.2092227225229806
$arr[0,12] = 0.1487714340921897
$arr[0,13] = 1.905393724496388
$arr[0,14] = 3.38450055242663
$arr[1,0] = 1
$arr[1,1] = 0.4350089752580857
$arr[1,2] = 0.1104118661123081
$arr[1,3] = 0
$arr[1,4] = 0.1121253597160941
$arr[1,5] = 0.387822817061874
$arr[1,6] = 0.002470289708980686
$arr[1,7] = 0
$arr[1,8] = 0.8963716224442386
$arr[1,9] = 0
$arr[1,10] = 0.2368935610889338
$arr[1,11] = 0.2066062862443445
$arr[1,12] = 0.1425545601854736
$arr[1,13] = 1.922617475217256
$arr[1,14] = 3.410770926026757
$arr[2,0] = 2
$arr[2,1] = 0.4156283381067283
$arr[2,2] = 0.1096859893203046
$arr[2,3] = 0
$arr[2,4] = 0.1123572962055626
$arr[2,5] = 0.3531389305168915
$arr[2,6] = 0.002471597993802521
$arr[2,7] = 0
$arr[2,8] = 0.9015012131756563
$arr[2,9] = 0
$arr[2,10] = 0.2190762765505951
$arr[2,11] = 0.2050987868208054
$arr[2,12] = 0.1387949932688777
$arr[2,13] = 1.933761047539985
$arr[2,14] = 3.428468009725052
$arr[3,0] = 3
$arr[3,1] = 0.407761959234449
$arr[3,2] = 0.1093885108958403
$arr[3,3] = 0
$arr[3,4] = 0.1124625768779364
$arr[3,5] = 0.3390132514313251
$arr[3,6] = 0.00247214801053246
$arr[3,7] = 0
$arr[3,8] = 0.9036859761152023
$arr[3,9] = 0
$arr[3,10] = 0.2118173505133001
$arr[3,11] = 0.204509446300122
$arr[3,12] = 0.137277558966133
$arr[3,13] = 1.938445086732903
$arr[3,14] = 3.436073944384631
$arr[4,0] = 4
$arr[4,1] = 0.4064576670137683
$arr[4,2] = 0.1093390139411667
$arr[4,3] = 0
$arr[4,4] = 0.1124807095097449
$arr[4,5] = 0.336668177824194
$arr[4,6] = 0.002472240361792266
$arr[4,7] = 0
$arr[4,8] = 0.9040544583428556
$arr[4,9] = 0
$arr[4,10] = 0.2106121316206782
$arr[4,11] = 0.2044130977775112
$arr[4,12] = 0.1370264770060494
$arr[4,13] = 1.939231503800967
$arr[4,14] = 3.437360720498503
$arr[5,0] = 5
$arr[5,1] = 0.415522121460242
$arr[5,2] = 0.1096819841946512
$arr[5,3] = 0
$arr[5,4] = 0.112358672440088
$arr[5,5] = 0.3529483938344953
$arr[5,6] = 0.002471605343457705
$arr[5,7] = 0
$arr[5,8] = 0.9015302953013737
$arr[5,9] = 0
$arr[5,10] = 0.2189783723788281
$arr[5,11] = 0.2050907375247846
$arr[5,12] = 0.1387744692506736
$arr[5,13] = 1.933823639245796
$arr[5,14] = 3.428568989530419
$arr[6,0] = 6
$arr[6,1] = 0.4557954710687682
$arr[6,2] = 0.1111806673481865
$arr[6,3] = 0
$arr[6,4] = 0.1119145909403709
$arr[6,5] = 0.4248636149813336
$arr[6,6] = 0.002468951209410298
$arr[6,7] = 0
$arr[6,8] = 0.8912204323822515
$arr[6,9] = 0
$arr[6,10] = 0.2559109682901664
$arr[6,11] = 0.2083000755116231
$arr[6,12] = 0.1466159736198094
$arr[6,13] = 1.91121461274502
$arr[6,14] = 3.393233473945202
$arr[7,0] = 7
$arr[7,1] = 0.5357191041066471
$arr[7,2] = 0.1140656363810209
$arr[7,3] = 0
$arr[7,4] = 0.1113816881046112
$arr[7,5] = 0.5661985755041457
$arr[7,6] = 0.002464275371056057
$arr[7,7] = 0
$arr[7,8] = 0.8739775801901217
$arr[7,9] = 0
$arr[7,10] = 0.3283593075703095
$arr[7,11] = 0.2153760279264318
$arr[7,12] = 0.1624454345336055
$arr[7,13] = 1.871381791131057
$arr[7,14] = 3.336365498198077
$arr[8,0] = 8
$arr[8,1] = 0.5949916297522577
$arr[8,2] = 0.1161518109552517
$arr[8,3] = 0
$arr[8,4] = 0.1111949484317645
$arr[8,5] = 0.6702781546542269
$arr[8,6] = 0.00246115938351985
$arr[8,7] = 0
$arr[8,8] = 0.8631192466017303
$arr[8,9] = 0
$arr[8,10] = 0.3815835927494788
$arr[8,11] = 0.2210483633700875
$arr[8,12] = 0.1743457897000908
$arr[8,13] = 1.844855194741289
$arr[8,14] = 3.302148781571944
$arr[9,0] = 9
$arr[9,1] = 0.6220712578676455
$arr[9,2] = 0.11709351057916
$arr[9,3] = 0
$arr[9,4] = 0.1111542071838691
$arr[9,5] = 0.7176906081379002
$arr[9,6] = 0.002459810527859002
$arr[9,7] = 0
$arr[9,8] = 0.8585718016773782
$arr[9,9] = 0
$arr[9,10] = 0.4057924332531115
$arr[9,11] = 0.2237310317277093
$arr[9,12] = 0.179817249520255
$arr[9,13] = 1.833380769016955
$arr[9,14] = 3.288223122266032
$arr[10,0] = 10
$arr[10,1] = 0.6323417714075106
$arr[10,2] = 0.1174490435593967
$arr[10,3] = 0
$arr[10,4] = 0.1111451139288349
$arr[10,5] = 0.7356546913071611
$arr[10,6] = 0.002459309569820739
$arr[10,7] = 0
$arr[10,8] = 0.8569061293934048
$arr[10,9] = 0
$arr[10,10] = 0.4149588238311708
$arr[10,11] = 0.2247615249433892
$arr[10,12] = 0.1818973621287014
$arr[10,13] = 1.829120845778242
$arr[10,14] = 3.283185460619194
$arr[11,0] = 11
$arr[11,1] = 0.6301291295993394
$arr[11,2] = 0.1173725209482441
$arr[11,3] = 0
$arr[11,4] = 0.1111467909493769
$arr[11,5] = 0.7317853510981394
$arr[11,6] = 0.002459417023776149
$arr[11,7] = 0
$arr[11,8] = 0.8572623562098585
$arr[11,9] = 0
$arr[11,10] = 0.4129847290965358
$arr[11,11] = 0.2245389406983804
$arr[11,12] = 0.1814490106689846
$arr[11,13] = 1.830034507956018
$arr[11,14] = 3.284259930751176
$arr[12,0] = 12
$arr[12,1] = 0.6229159014595211
$arr[12,2] = 0.11712278199672
$arr[12,3] = 0
$arr[12,4] = 0.1111533322504812
$arr[12,5] = 0.7191683204515869
$arr[12,6] = 0.002459769117274402
$arr[12,7] = 0
$arr[12,8] = 0.8584336367331922
$arr[12,9] = 0
$arr[12,10] = 0.4065465801718915
$arr[12,11] = 0.2238155184064539
$arr[12,12] = 0.1799882183118555
$arr[12,13] = 1.833028595344409
$arr[12,14] = 3.287803947955553
$arr[13,0] = 13
$arr[13,1] = 0.6184996606762354
$arr[13,2] = 0.1169696700650036
$arr[13,3] = 0
$arr[13,4] = 0.1111581632689926
$arr[13,5] = 0.7114413442032514
$arr[13,6] = 0.002459986062578462
$arr[13,7] = 0
$arr[13,8] = 0.8591584168909598
$arr[13,9] = 0
$arr[13,10] = 0.4026028860936037
$arr[13,11] = 0.2233743034379643
$arr[13,12] = 0.1790945031570175
$arr[13,13] = 1.834873654071892
$arr[13,14] = 3.290005450693073
$arr[14,0] = 14
$arr[14,1] = 0.5932241943728513
$arr[14,2] = 0.1160901201015605
$arr[14,3] = 0
$arr[14,4] = 0.111198498620638
$arr[14,5] = 0.6671810134426437
$arr[14,6] = 0.002461248911340899
$arr[14,7] = 0
$arr[14,8] = 0.86342431943077
$arr[14,9] = 0
$arr[14,10] = 0.3800013803185038
$arr[14,11] = 0.2208750958601371
$arr[14,12] = 0.1739893705944695
$arr[14,13] = 1.845616998541832
$arr[14,14] = 3.303091842229023
$arr[15,0] = 15
$arr[15,1] = 0.5777478132291094
$arr[15,2] = 0.1155486600153779
$arr[15,3] = 0
$arr[15,4] = 0.1112345497169329
$arr[15,5] = 0.6400460337125793
$arr[15,6] = 0.00246204117333308
$arr[15,7] = 0
$arr[15,8] = 0.8661417069119501
$arr[15,9] = 0
$arr[15,10] = 0.3661349120720274
$arr[15,11] = 0.2193680527136905
$arr[15,12] = 0.1708722740778086
$arr[15,13] = 1.852359463363253
$arr[15,14] = 3.31153981178916
$arr[16,0] = 16
$arr[16,1] = 0.5688571911409213
$arr[16,2] = 0.1152365392894126
$arr[16,3] = 0
$arr[16,4] = 0.1112594477205615
$arr[16,5] = 0.6244449056556647
$arr[16,6] = 0.002462503322948637
$arr[16,7] = 0
$arr[16,8] = 0.8677415812832621
$arr[16,9] = 0
$arr[16,10] = 0.3581590204490226
$arr[16,11] = 0.2185108748916775
$arr[16,12] = 0.169084863528326
$arr[16,13] = 1.856293336047894
$arr[16,14] = 3.316553202593781
$arr[17,0] = 17
$arr[17,1] = 0.5658488860737805
$arr[17,2] = 0.1151307430599502
$arr[17,3] = 0
$arr[17,4] = 0.1112685933304025
$arr[17,5] = 0.619163680173358
$arr[17,6] = 0.002462660909914741
$arr[17,7] = 0
$arr[17,8] = 0.8682896107021953
$arr[17,9] = 0
$arr[17,10] = 0.3554584907234357
$arr[17,11] = 0.2182223063609143
$arr[17,12] = 0.1684806192264432
$arr[17,13] = 1.857634859886545
$arr[17,14] = 3.318277162141356
$arr[18,0] = 18
$arr[18,1] = 0.5793941685874131
$arr[18,2] = 0.1156063706112604
$arr[18,3] = 0
$arr[18,4] = 0.1112302813894832
$arr[18,5] = 0.642933953830422
$arr[18,6] = 0.002461956167141936
$arr[18,7] = 0
$arr[18,8] = 0.8658486170836071
$arr[18,9] = 0
$arr[18,10] = 0.3676110525559011
$arr[18,11] = 0.2195274836090277
$arr[18,12] = 0.1712035301558998
$arr[18,13] = 1.851635943136252
$arr[18,14] = 3.310624537920432
$arr[19,0] = 19
$arr[19,1] = 0.6250341712849945
$arr[19,2] = 0.1171961655057743
$arr[19,3] = 0
$arr[19,4] = 0.1111512391722229
$arr[19,5] = 0.7228739723491628
$arr[19,6] = 0.002459665432428091
$arr[19,7] = 0
$arr[19,8] = 0.8580880742805874
$arr[19,9] = 0
$arr[19,10] = 0.4084376510847108
$arr[19,11] = 0.2240276088311219
$arr[19,12] = 0.1804170672969647
$arr[19,13] = 1.832146847424536
$arr[19,14] = 3.286756588637445
$arr[20,0] = 20
$arr[20,1] = 0.6549557682792795
$arr[20,2] = 0.1182289541349633
$arr[20,3] = 0
$arr[20,4] = 0.1111364928346035
$arr[20,5] = 0.7751780083420101
$arr[20,6] = 0.00245822554854768
$arr[20,7] = 0
$arr[20,8] = 0.8533445262265751
$arr[20,9] = 0
$arr[20,10] = 0.4351142874242271
$arr[20,11] = 0.2270539065111024
$arr[20,12] = 0.1864863158866967
$arr[20,13] = 1.819906152930223
$arr[20,14] = 3.272531199042106
$arr[21,0] = 21
$arr[21,1] = 0.638977741527043
$arr[21,2] = 0.1176783113855535
$arr[21,3] = 0
$arr[21,4] = 0.1111409931620244
$arr[21,5] = 0.7472568307830727
$arr[21,6] = 0.002458988818240086
$arr[21,7] = 0
$arr[21,8] = 0.8558462076301296
$arr[21,9] = 0
$arr[21,10] = 0.4208771763198911
$arr[21,11] = 0.2254309462738888
$arr[21,12] = 0.1832427284129423
$arr[21,13] = 1.826393816455234
$arr[21,14] = 3.279997897995145
$arr[22,0] = 22
$arr[22,1] = 0.5786498296134539
$arr[22,2] = 0.1155802822284144
$arr[22,3] = 0
$arr[22,4] = 0.1112321981043163
$arr[22,5] = 0.6416283278902171
$arr[22,6] = 0.002461994577659151
$arr[22,7] = 0
$arr[22,8] = 0.8659810058340121
$arr[22,9] = 0
$arr[22,10] = 0.3669437014719108
$arr[22,11] = 0.2194553761094937
$arr[22,12] = 0.1710537547890638
$arr[22,13] = 1.851962867395528
$arr[22,14] = 3.311037845616909
$arr[23,0] = 23
$arr[23,1] = 0.5139987806832096
$arr[23,2] = 0.1132910054899412
$arr[23,3] = 0
$arr[23,4] = 0.1114898011664085
$arr[23,5] = 0.5279251897347166
$arr[23,6] = 0.002465484006342347
$arr[23,7] = 0
$arr[23,8] = 0.8783241236501809
$arr[23,9] = 0
$arr[23,10] = 0.3087594336315931
$arr[23,11] = 0.2133783761595254
$arr[23,12] = 0.1581152700910806
$arr[23,13] = 1.881676218955295
$arr[23,14] = 3.350420697966271

$ws.Range("A2:O25").Value2 = $arr
